$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 172
$ws.Range("I11").Value = 172
$ws.Range("K11").Value = 172
$ws.Range("M11").Value = -32
$ws.Range("H64").Value = 4800
$ws.Range("H67").Value = 4800
$ws.Range("H92").Value = 643.17645
$ws.Range("I92").Value = 698.8333
$ws.Range("J92").Value = 509.6
$ws.Range("K92").Value = 698.8333
$ws.Range("L92").Value = 509.6
$ws.Range("M92").Value = 549.1667
$ws.Range("N92").Value = -3005.6
$ws.Range("H104").Value = 591.5
$ws.Range("I104").Value = 591.5
$ws.Range("K104").Value = 1774.5
$ws.Range("M104").Value = -27.5
$ws.Range("H113").Value = 5806.7144
$ws.Range("I113").Value = 6719.7
$ws.Range("J113").Value = 3524.25
$ws.Range("K113").Value = 6719.7
$ws.Range("L113").Value = 3524.25
$ws.Range("M113").Value = -3465.7
$ws.Range("N113").Value = -10032.25
$ws.Range("H137").Value = 13268.527
$ws.Range("I137").Value = 5214.4
$ws.Range("J137").Value = 23336.188
$ws.Range("K137").Value = 15643.2
$ws.Range("L137").Value = 70008.564
$ws.Range("M137").Value = -13093.2
$ws.Range("N137").Value = -75108.564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4387.6113
$ws.Range("I61").Value = 1509.091
$ws.Range("J61").Value = 8911
$ws.Range("K61").Value = 1509.091
$ws.Range("L61").Value = 8911
$ws.Range("M61").Value = -1297.091
$ws.Range("N61").Value = -9335
$ws.Range("H74").Value = 4507.222
$ws.Range("I74").Value = 2865.7058
$ws.Range("K74").Value = 2865.7058
$ws.Range("M74").Value = -1991.7058
$ws.Range("H77").Value = 4507.222
$ws.Range("I77").Value = 2865.7058
$ws.Range("K77").Value = 14328.529
$ws.Range("M77").Value = -9960.529
$ws.Range("H122").Value = 1840.8334
$ws.Range("I122").Value = 1879.16
$ws.Range("J122").Value = 1753.7273
$ws.Range("K122").Value = 5637.48
$ws.Range("L122").Value = 5261.1819
$ws.Range("M122").Value = -3187.48
$ws.Range("N122").Value = -10161.1819
$ws.Range("H136").Value = 4387.6113
$ws.Range("I136").Value = 1509.091
$ws.Range("J136").Value = 8911
$ws.Range("K136").Value = 4527.272999999999
$ws.Range("L136").Value = 26733
$ws.Range("M136").Value = -1977.272999999999
$ws.Range("N136").Value = -31833
$ws.Range("H138").Value = 44474.5
$ws.Range("J138").Value = 44474.5
$ws.Range("L138").Value = 44474.5
$ws.Range("N138").Value = -54754.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 16161
$ws.Range("I29").Value = 16161
$ws.Range("K29").Value = 16161
$ws.Range("M29").Value = -15872
$ws.Range("H82").Value = 17015.143
$ws.Range("J82").Value = 21462.25
$ws.Range("L82").Value = 21462.25
$ws.Range("N82").Value = -22228.25
$ws.Range("H85").Value = 17015.143
$ws.Range("J85").Value = 21462.25
$ws.Range("L85").Value = 21462.25
$ws.Range("N85").Value = -24114.25
$ws.Range("H99").Value = 2153.9092
$ws.Range("I99").Value = 2062.375
$ws.Range("K99").Value = 2062.375
$ws.Range("M99").Value = -564.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1603.4333
$ws.Range("J31").Value = 3220.3333
$ws.Range("L31").Value = 3220.3333
$ws.Range("N31").Value = -3810.3333
$ws.Range("H34").Value = 1603.4333
$ws.Range("J34").Value = 3220.3333
$ws.Range("L34").Value = 3220.3333
$ws.Range("N34").Value = -3624.3333
$ws.Range("H58").Value = 14641.571
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 14641.571
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 14641.571
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -15047.571
$ws.Range("H99").Value = 7990.96
$ws.Range("I99").Value = 5019.174
$ws.Range("K99").Value = 5019.174
$ws.Range("M99").Value = -3521.174
$ws.Range("H107").Value = 832.8461
$ws.Range("J107").Value = 939.5714
$ws.Range("L107").Value = 939.5714
$ws.Range("N107").Value = -4779.5714
$ws.Range("H126").Value = 7990.96
$ws.Range("I126").Value = 5019.174
$ws.Range("K126").Value = 15057.522
$ws.Range("M126").Value = -12587.522
$ws.Range("H136").Value = 14641.571
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 14641.571
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 43924.713
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -49024.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 64365044
$ws.Range("I4").Value = 68891120
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 206673360
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -206673248
$ws.Range("N4").Value = -3000224
$ws.Range("H6").Value = 155.5
$ws.Range("I6").Value = 152.25
$ws.Range("J6").Value = 175
$ws.Range("K6").Value = 456.75
$ws.Range("L6").Value = 525
$ws.Range("M6").Value = -343.75
$ws.Range("N6").Value = -751
$ws.Range("H14").Value = 2433.5454
$ws.Range("I14").Value = 2433.5454
$ws.Range("K14").Value = 7300.6362
$ws.Range("M14").Value = -7127.6362
$ws.Range("H60").Value = 2249
$ws.Range("J60").Value = 2498.5
$ws.Range("L60").Value = 7495.5
$ws.Range("N60").Value = -7997.5
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H140").Value = 3232.5715
$ws.Range("I140").Value = 1907
$ws.Range("K140").Value = 5721
$ws.Range("M140").Value = -541

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 119.21429
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 207.8
$ws.Range("K2").Value = 70
$ws.Range("L2").Value = 207.8
$ws.Range("M2").Value = 43
$ws.Range("N2").Value = -433.8
$ws.Range("H70").Value = 7073.4116
$ws.Range("I70").Value = 5512
$ws.Range("K70").Value = 5512
$ws.Range("M70").Value = -5242
$ws.Range("H73").Value = 7073.4116
$ws.Range("I73").Value = 5512
$ws.Range("K73").Value = 5512
$ws.Range("M73").Value = -4576
$ws.Range("H126").Value = 3209.5557
$ws.Range("J126").Value = 3733.3333
$ws.Range("L126").Value = 11199.9999
$ws.Range("N126").Value = -16139.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2637.1052
$ws.Range("J22").Value = 3483.4546
$ws.Range("L22").Value = 3483.4546
$ws.Range("N22").Value = -4073.4546
$ws.Range("H24").Value = 24749.5
$ws.Range("J24").Value = 24749.5
$ws.Range("L24").Value = 24749.5
$ws.Range("N24").Value = -25435.5
$ws.Range("H27").Value = 2637.1052
$ws.Range("J27").Value = 3483.4546
$ws.Range("L27").Value = 3483.4546
$ws.Range("N27").Value = -3697.4546
$ws.Range("H43").Value = 437191
$ws.Range("I43").Value = 321257
$ws.Range("J43").Value = 553125
$ws.Range("K43").Value = 321257
$ws.Range("L43").Value = 553125
$ws.Range("M43").Value = -321064
$ws.Range("N43").Value = -553511

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 15007
$ws.Range("J18").Value = 15007
$ws.Range("L18").Value = 15007
$ws.Range("N18").Value = -15353
$ws.Range("H22").Value = 9163.571
$ws.Range("I22").Value = 5306.5
$ws.Range("J22").Value = 10706.4
$ws.Range("K22").Value = 5306.5
$ws.Range("L22").Value = 10706.4
$ws.Range("M22").Value = -5013.5
$ws.Range("N22").Value = -11292.4
$ws.Range("H31").Value = 13346
$ws.Range("I31").Value = 20000
$ws.Range("J31").Value = 10019
$ws.Range("K31").Value = 20000
$ws.Range("L31").Value = 10019
$ws.Range("M31").Value = -19652
$ws.Range("N31").Value = -10715
$ws.Range("H70").Value = 49280.43
$ws.Range("I70").Value = 43497.5
$ws.Range("J70").Value = 51593.6
$ws.Range("K70").Value = 43497.5
$ws.Range("L70").Value = 51593.6
$ws.Range("M70").Value = -43182.5
$ws.Range("N70").Value = -52223.6
$ws.Range("H73").Value = 49280.43
$ws.Range("I73").Value = 43497.5
$ws.Range("J73").Value = 51593.6
$ws.Range("K73").Value = 43497.5
$ws.Range("L73").Value = 51593.6
$ws.Range("M73").Value = -42405.5
$ws.Range("N73").Value = -53777.6
$ws.Range("H81").Value = 2969.625
$ws.Range("J81").Value = 2979.8
$ws.Range("L81").Value = 5959.6
$ws.Range("N81").Value = -8081.6
$ws.Range("H84").Value = 2969.625
$ws.Range("J84").Value = 2979.8
$ws.Range("L84").Value = 29798
$ws.Range("N84").Value = -40406
$ws.Range("H105").Value = 40257.5
$ws.Range("J105").Value = 40257.5
$ws.Range("L105").Value = 40257.5
$ws.Range("N105").Value = -47245.5
$ws.Range("H107").Value = 1363.7778
$ws.Range("J107").Value = 1645.125
$ws.Range("L107").Value = 4935.375
$ws.Range("N107").Value = -8775.375
$ws.Range("H113").Value = 399.22223
$ws.Range("I113").Value = 375.83334
$ws.Range("J113").Value = 446
$ws.Range("K113").Value = 1127.50002
$ws.Range("L113").Value = 1338
$ws.Range("M113").Value = 1042.49998
$ws.Range("N113").Value = -5678
$ws.Range("H122").Value = 3268.9565
$ws.Range("I122").Value = 3399.7058
$ws.Range("K122").Value = 10199.1174
$ws.Range("M122").Value = -7749.117400000001
$ws.Range("H126").Value = 5575.032
$ws.Range("I126").Value = 5522.12
$ws.Range("K126").Value = 16566.36
$ws.Range("M126").Value = -14096.36
